$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellAddr, $value) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '58.269.06'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '2.597.34'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '520.56'
$ws.Range('E5').Value = '  +0.79%  '
Set-TextValue 'D6' '144.62'
$ws.Range('E6').Value = '  +2.19%  '
Set-TextValue 'D7' '0.997'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('D9').Value = '2.615.85'
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').Value = '3.054.57'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '58.267.71'
$ws.Range('E15').Value = '  +0.38%  '
Set-TextValue 'D16' '20.57'
$ws.Range('E16').Value = '  -0.14%  '
Set-TextValue 'D17' '0.0000135'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '2.607.86'
$ws.Range('E18').Value = '  -0.24%  '
Set-TextValue 'D19' '341.11'
$ws.Range('E19').Value = '  +2.06%  '
Set-TextValue 'D20' '4.35'
$ws.Range('E20').Value = '  -0.66%  '
Set-TextValue 'D21' '10.33'
$ws.Range('E21').Value = '  +0.09%  '
Set-TextValue 'D22' '6.39'
$ws.Range('E22').Value = '  +2.40%  '
Set-TextValue 'D23' '0.999'
$ws.Range('E23').Value = '  +0.07%  '
Set-TextValue 'D24' '66.35'
$ws.Range('E24').Value = '  +3.49%  '
$ws.Range('E25').Value = '  +0.17%  '
Set-TextValue 'D26' '0.404'
$ws.Range('E26').Value = '  -4.29%  '
Set-TextValue 'D27' '0.997'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '2.713.37'
$ws.Range('E28').Value = '  -1.08%  '
Set-TextValue 'D29' '7.05'
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('D30').Value = '0.0₃0755'
$ws.Range('E30').Value = '  -3.06%  '
$ws.Range('E31').Value = '  -0.10%  '
Set-TextValue 'D32' '6.26'
$ws.Range('E32').Value = '  -5.04%  '
Set-TextValue 'D33' '1.59'
$ws.Range('E33').Value = '  +0.77%  '
Set-TextValue 'D34' '18.83'
$ws.Range('E34').Value = '  +0.97%  '
Set-TextValue 'D35' '149.85'
$ws.Range('E35').Value = '  -1.15%  '
Set-TextValue 'D36' '4.05'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').Value = '  -1.32%  '
Set-TextValue 'D38' '0.880'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D39' '0.846'
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D40' '1.47'
$ws.Range('E40').Value = '  +3.18%  '
Set-TextValue 'D41' '36.11'
$ws.Range('E41').Value = '  -1.38%  '
$ws.Range('E42').Value = '  -0.55%  '
Set-TextValue 'D43' '0.996'
$ws.Range('E43').Value = '  -0.52%  '
Set-TextValue 'D44' '275.09'
$ws.Range('E44').Value = '  +2.97%  '
$ws.Range('E45').Value = '  -0.09%  '
Set-TextValue 'D46' '0.0957'
$ws.Range('E46').Value = '  -0.63%  '
Set-TextValue 'D48' '18.95'
$ws.Range('E48').Value = '  -0.65%  '
Set-TextValue 'D49' '0.0524'
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D50' '4.73'
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D51' '19.12'
$ws.Range('E51').Value = '  +5.49%  '
